$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that can be set directly as text (percentages, coin names, URLs)
$directUpdates = @{
    'E2' = '  +0.02%  '
    'E3' = '  +0.20%  '
    'E4' = '  +0.05%  '
    'E5' = '  +0.46%  '
    'E6' = '  -2.64%  '
    'E7' = '  +5.84%  '
    'E8' = '  -0.06%  '
    'E9' = '  +6.43%  '
    'E10' = '  +5.48%  '
    'E11' = '  +1.33%  '
    'E12' = '  +8.29%  '
    'B13' = 'TRON'
    'C13' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'E13' = '  -0.02%  '
    'B14' = 'Chainlink'
    'C14' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E14' = '  +7.38%  '
    'E15' = '  +0.29%  '
    'E16' = '  +41.40%  '
    'E17' = '  +1.02%  '
    'E18' = '  +5.57%  '
    'E19' = '  +5.63%  '
    'E20' = '  +0.30%  '
    'E21' = '  +41.26%  '
    'E22' = '  +9.00%  '
    'E23' = '  +0.37%  '
    'E24' = '  +1.42%  '
    'E25' = '  +2.82%  '
    'E26' = '  +10.43%  '
    'E27' = '  +6.85%  '
    'E28' = '  -0.30%  '
    'E29' = '  -3.59%  '
    'E30' = '  +0.29%  '
    'E31' = '  +5.68%  '
    'E32' = '  -0.67%  '
    'E33' = '  -0.73%  '
    'E34' = '  -2.38%  '
    'E35' = '  -0.07%  '
    'E36' = '  +2.35%  '
    'E37' = '  +3.54%  '
    'E38' = '  +0.06%  '
    'E39' = '  -0.01%  '
    'E40' = '  +7.26%  '
    'E41' = '  +2.23%  '
    'E42' = '  -2.59%  '
    'E43' = '  +1.62%  '
    'E44' = '  +6.99%  '
    'E45' = '  +0.16%  '
    'E46' = '  +13.08%  '
    'E47' = '  -1.09%  '
    'E48' = '  +4.70%  '
    'E49' = '  +22.19%  '
    'E50' = '  +10.85%  '
    'E51' = '  +0.72%  '
}

# Cells that look numeric (single-decimal prices) - need apostrophe prefix
# to force text storage and prevent Excel from parsing/rounding them as numbers,
# then Style is reset to Normal so no stray quote-prefix formatting is applied.
$textUpdates = @{
    'D2' = '61.907.15'
    'D3' = '3.423.85'
    'D5' = '408.24'
    'D6' = '128.18'
    'D7' = '0.629'
    'D11' = '42.52'
    'D12' = '9.11'
    'D13' = '0.141'
    'D14' = '21.43'
    'D15' = '3.964.10'
    'D16' = '0.0000206'
    'D17' = '3.435.55'
    'D18' = '12.34'
    'D19' = '1.08'
    'D20' = '61.954.35'
    'D21' = '442.47'
    'D22' = '91.52'
    'D24' = '12.96'
    'D26' = '32.74'
    'D27' = '8.68'
    'D29' = '7.73'
    'D30' = '2.72'
    'D31' = '11.98'
    'D34' = '42.90'
    'D35' = '0.999'
    'D37' = '53.28'
    'D41' = '0.320'
    'D43' = '141.72'
    'D44' = '4.23'
    'D46' = '2.50'
    'D47' = '16.57'
    'D48' = '22.26'
    'D50' = '2.14'
    'D51' = '3.772.96'
}

foreach ($cell in $directUpdates.Keys) {
    $ws.Range($cell).Value = $directUpdates[$cell]
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = "'" + $textUpdates[$cell]
    $ws.Range($cell).Style = "Normal"
}

Write-Output "Applied $($directUpdates.Count + $textUpdates.Count) cell updates"